$d = $word.ActiveDocument

$find = "Сазвежђе Боотес цонстеллатион 2022: 14-23. мај, 13-22. јун, 12-21. јул"
$replace = "Сазвежђе Боотес цонстеллатион током 2022. године посматрамо 14-23. мај, 13-22. јун, 12-21. јул"

$r = $d.Content
$r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
